$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: reorder the "Value n" placeholder columns and append a new ---
# --- trailing "Value 7" column (L), growing the sheet from A1:K30 to A1:L30. ---
$ws.Range("E1").Value = "Value 2"
$ws.Range("F1").Value = "Value 3"
$ws.Range("G1").Value = "Value n"
$ws.Range("H1").Value = "Value 1"
$ws.Range("I1").Value = "Value 4"
$ws.Range("J1").Value = "Value 5"
$ws.Range("K1").Value = "Value 6"
$ws.Range("L1").Value = "Value 7"

# New column L needs the same header styling as the rest of row 1 (bold, centered,
# thin box border) since it did not exist before.
$ws.Range("L1").Font.Bold = $true
$ws.Range("L1").HorizontalAlignment = -4108
$ws.Range("L1").VerticalAlignment = -4160
$ws.Range("L1").Borders.LineStyle = 1

# --- Column "Value" (D) content edits, plus a few column E/F companions ---
$ws.Range("E3").Value = "Experimental"
$ws.Range("D5").Value = "Duke University DCE-MRI of breast cancer patients"
$ws.Range("D6").Value = "Retrospective collection of MRI from 922 biopsy-confirmed invasive breast cancer patients."
$ws.Range("D11").Value = "Breast MRI is a medical image modality used to assess the extent of disease in breast cancer patients. Recent studies show that MRI has a potential in prognosis of patients’ short and long-term outcomes as well as predicting pathological and genomic features of the tumors. However, large, well annotated datasets are needed to make further progress in the field. We share such a dataset here."
$ws.Range("D12").Value = "This dataset is a single-institutional, retrospective collection of 922 biopsy-confirmed invasive breast cancer patients, over a decade, specifically pre-operative dynamic contrast enhanced (DCE)-MRI that were downloaded from PACS systems and de-identified for The Cancer Imaging Archive (TCIA) release. These include axial breast MRI images acquired by 1.5T or 3T scanners in the prone positions. The following MRI sequences are shared in DICOM format: a non-fat saturated T1-weighted sequence, a fat-saturated gradient echo T1-weighted pre-contrast sequence, and mostly three to four post-contrast sequences."
$ws.Range("D13").Value = "Data collected for subsequent analysis."
$ws.Range("D15").Value = "Imaging"
$ws.Range("D16").Value = "MRI"
$ws.Range("D22").Value = "PrincipalInvestigator"
$ws.Range("E22").Value = "CorrespondingAuthor"
$ws.Range("D24").Value = "Protocol for dataset"
$ws.Range("E24").Value = "Paper for dataset"
$ws.Range("F24").Value = "Not Defined"
$ws.Range("D25").Value = "HasProtocol"
$ws.Range("E25").Value = "IsDescribedBy"
$ws.Range("F25").Value = "Not Defined"
$ws.Range("D26").Value = "Not Defined"
$ws.Range("D27").Value = "Not Defined"
$ws.Range("D29").Value = 2
$ws.Range("D30").Value = 4

# --- Cells that no longer hold a value ---
$ws.Range("D17").Value = $null
$ws.Range("F22").Value = $null
$ws.Range("G22").Value = $null
$ws.Range("H22").Value = $null
$ws.Range("I22").Value = $null
$ws.Range("J22").Value = $null
